# Insert a new row at position 797 (pushes existing rows 797-838 down to 798-839)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(797).Insert()

# Set the date as literal text (avoid Excel auto-converting "2026/02/08" into a date serial value)
$ws.Range("A797").NumberFormat = "@"
$ws.Range("A797").Value = "2026/02/08"
$ws.Range("A797").ClearFormats()

$ws.Range("B797").Value = "日"
$ws.Range("C797").Value = 0
$ws.Range("D797").Value = 201
